$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G2").Value = 5.416159
$ws.Range("H2").Value = 16.248477
$ws.Range("I2").Value = 0.07167636930710344
$ws.Range("J2").Value = 0.07167636930710343
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.1837056666666667
$ws.Range("N2").Value = 0.5511170000000001
$ws.Range("O2").Value = 0.7269991860920679
$ws.Range("P2").Value = 0.7269991860920678
$ws.Range("Q2").Value = 0.9949790998676669
$ws.Range("R2").Value = 8.954811898809002
$ws.Range("S2").Value = 0.05210866214829868
$ws.Range("T2").Value = 0.05210866214829867
$ws.Range("G3").Value = 5.416159
$ws.Range("H3").Value = 16.248477
$ws.Range("I3").Value = 0.07167636930710344
$ws.Range("J3").Value = 0.07167636930710343
$ws.Range("O3").Value = 0.2534828531892131
$ws.Range("P3").Value = 0.2534828531892131
$ws.Range("Q3").Value = 0.3469194270406666
$ws.Range("R3").Value = 3.122274843366
$ws.Range("S3").Value = 0.01816873059820832
$ws.Range("T3").Value = 0.01816873059820832
$ws.Range("G4").Value = 5.416159
$ws.Range("H4").Value = 16.248477
$ws.Range("I4").Value = 0.07167636930710344
$ws.Range("J4").Value = 0.07167636930710343
$ws.Range("O4").Value = 0.01951796071871896
$ws.Range("P4").Value = 0.01951796071871896
$ws.Range("Q4").Value = 0.026712496188
$ws.Range("R4").Value = 0.240412465692
$ws.Range("S4").Value = 0.001398976560596438
$ws.Range("T4").Value = 0.001398976560596438
$ws.Range("G5").Value = 61.15258266666666
$ws.Range("I5").Value = 0.809281097415931
$ws.Range("J5").Value = 0.8092810974159309
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.1837056666666667
$ws.Range("N5").Value = 0.5511170000000001
$ws.Range("O5").Value = 0.7269991860920679
$ws.Range("P5").Value = 0.7269991860920678
$ws.Range("Q5").Value = 11.23407596716844
$ws.Range("R5").Value = 101.106683704516
$ws.Range("S5").Value = 0.5883466991410774
$ws.Range("T5").Value = 0.5883466991410772
$ws.Range("G6").Value = 61.15258266666666
$ws.Range("I6").Value = 0.809281097415931
$ws.Range("J6").Value = 0.8092810974159309
$ws.Range("O6").Value = 0.2534828531892131
$ws.Range("P6").Value = 0.2534828531892131
$ws.Range("S6").Value = 0.2051388816050877
$ws.Range("T6").Value = 0.2051388816050877
$ws.Range("G7").Value = 61.15258266666666
$ws.Range("I7").Value = 0.809281097415931
$ws.Range("J7").Value = 0.8092810974159309
$ws.Range("O7").Value = 0.01951796071871896
$ws.Range("P7").Value = 0.01951796071871896
$ws.Range("S7").Value = 0.01579551666976591
$ws.Range("T7").Value = 0.01579551666976591
$ws.Range("I8").Value = 0.1190425332769656
$ws.Range("J8").Value = 0.1190425332769656
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.1837056666666667
$ws.Range("N8").Value = 0.5511170000000001
$ws.Range("O8").Value = 0.7269991860920679
$ws.Range("P8").Value = 0.7269991860920678
$ws.Range("Q8").Value = 1.652494870358111
$ws.Range("R8").Value = 14.872453833223
$ws.Range("S8").Value = 0.0865438248026919
$ws.Range("T8").Value = 0.08654382480269188
$ws.Range("I9").Value = 0.1190425332769656
$ws.Range("J9").Value = 0.1190425332769656
$ws.Range("O9").Value = 0.2534828531892131
$ws.Range("P9").Value = 0.2534828531892131
$ws.Range("S9").Value = 0.03017524098591708
$ws.Range("T9").Value = 0.03017524098591708
$ws.Range("I10").Value = 0.1190425332769656
$ws.Range("J10").Value = 0.1190425332769656
$ws.Range("O10").Value = 0.01951796071871896
$ws.Range("P10").Value = 0.01951796071871896
$ws.Range("S10").Value = 0.002323467488356609
$ws.Range("T10").Value = 0.002323467488356608
